# Updating index tab templates
# Add new gradient/alignment columns to the "Custom Linear Assay Index " sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Custom Linear Assay Index ")

# --- Insert 4 new columns before the existing "XL coordinates" column (D) ---
# This shifts the previous D:G (XL/YL/XR/YR coordinates) to H:K
$ws.Range("D1:G1").EntireColumn.Insert()

# Populate headers for the newly inserted columns
$ws.Range("D1").Value = "Low gradient"
$ws.Range("E1").Value = "High gradient"
$ws.Range("F1").Value = "Distance between"
$ws.Range("G1").Value = "Gradient slope (value per per cm)"

# --- Append two new trailing columns (L, M) ---
$ws.Range("L1").Value = "Alignment distance"
$ws.Range("M1").Value = "Change orientation?"

# Match header formatting (centered) used by all the other header cells
$ws.Range("L1:M1").HorizontalAlignment = -4108

# --- Approximate the original author's column widths as closely as possible ---
# (ColumnWidth only supports coarse 1/6-character increments in this environment,
# so we pick the inputs whose rounded, persisted width is nearest the target.)
$ws.Columns("D:D").ColumnWidth = 10.998697916666666
$ws.Columns("E:E").ColumnWidth = 11.330729166666666
$ws.Columns("F:F").ColumnWidth = 14.998697916666666
$ws.Columns("G:G").ColumnWidth = 28.330729166666668
$ws.Columns("L:L").ColumnWidth = 16.330729166666668
$ws.Columns("M:M").ColumnWidth = 16.998697916666668

# Restore the selected cell shown in the saved workbook
$ws.Activate() | Out-Null
$ws.Range("D11").Select() | Out-Null
